$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark E41:E43 as "Done" (previously blank)
$ws.Range("E41").Value = "Done"
$ws.Range("E42").Value = "Done"
$ws.Range("E43").Value = "Done"

# New section header row 50 "April 30, 2018" (styled like row 40, merged A50:I50)
$ws.Range("A40:I40").Copy()
$ws.Range("A50:I50").PasteSpecial(-4122)
$ws.Range("A50").Value = "April 30, 2018"
$ws.Range("A50:I50").Merge()

# New data rows 51-53 (Dynamic Programming problems picked by I-No)
$ws.Range("A49:D49").Copy()
$ws.Range("A51:D51").PasteSpecial(-4122)
$ws.Range("A52:D52").PasteSpecial(-4122)
$ws.Range("A53:D53").PasteSpecial(-4122)
$ws.Range("F49:G49").Copy()
$ws.Range("F51:G51").PasteSpecial(-4122)
$ws.Range("F52:G52").PasteSpecial(-4122)
$ws.Range("F53:G53").PasteSpecial(-4122)

$ws.Range("A51").Value = 53
$ws.Range("B51").Value = "Maximum Subarray"
$ws.Range("C51").Value = "Dynamic Programming"
$ws.Range("D51").Value = "I-No"
$ws.Range("F51").Value = "Easy"
$ws.Range("G51").Value = "Python"

$ws.Range("A52").Value = 70
$ws.Range("B52").Value = "Climbing Stairs"
$ws.Range("C52").Value = "Dynamic Programming"
$ws.Range("D52").Value = "I-No"
$ws.Range("F52").Value = "Easy"
$ws.Range("G52").Value = "Python"

$ws.Range("A53").Value = 338
$ws.Range("B53").Value = "Counting Bits"
$ws.Range("C53").Value = "Dynamic Programming"
$ws.Range("D53").Value = "I-No"
$ws.Range("F53").Value = "Easy"
$ws.Range("G53").Value = "Python"

# Update view selection to match where the author ended up editing
[void]$ws.Range("E53").Select()
